$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before BI (27/04/2020 data), shifting BI..BL to BJ..BM.
# This mirrors adding a new "01/05/2020" reporting day into the wide table,
# which chronologically belongs right after "26/04/2020" (column BH) in the
# original (non-alphabetic) column ordering used by this sheet.
$ws.Range("BI1").EntireColumn.Insert()

# Set the new header text as a literal string "01/05/2020" (not an Excel date
# serial). Force text format first so Excel's locale-aware date parser does
# not reinterpret "01/05/2020" (a valid M/D/Y date) as a date value, then
# restore the default "Normal" style so no stray number-format is left on
# the cell.
$ws.Range("BI1").NumberFormat = "@"
$ws.Range("BI1").Value2 = "01/05/2020"
$ws.Range("BI1").Style = "Normal"

# Final values (after the column insert/shift) for columns BI..BM, rows 2-34.
# $null means the cell must be blank.
$rowData = @(
    ,@(2, $null, $null, 1, $null, 1)
    ,@(3, $null, 82, 85, 77, 75)
    ,@(4, $null, $null, $null, $null, $null)
    ,@(5, $null, $null, 2, 2, 2)
    ,@(6, $null, 69, 23, 37, 27)
    ,@(7, $null, 9, 11, 12, 7)
    ,@(8, $null, $null, 3, $null, 4)
    ,@(9, $null, 1, 2, 3, 3)
    ,@(10, $null, $null, $null, $null, $null)
    ,@(11, $null, 28, 35, 39, 37)
    ,@(12, $null, 8, 7, 4, 9)
    ,@(13, $null, 3, $null, $null, 5)
    ,@(14, $null, 17, 9, 8, 12)
    ,@(15, $null, 5, 3, 2, 1)
    ,@(16, $null, 21, 25, 22, 44)
    ,@(17, $null, 26, 8, 20, 16)
    ,@(18, $null, $null, 1, 1, $null)
    ,@(19, $null, 15, 21, 25, 16)
    ,@(20, $null, 33, 26, 26, 24)
    ,@(21, $null, $null, $null, $null, $null)
    ,@(22, $null, $null, $null, $null, $null)
    ,@(23, $null, $null, $null, $null, $null)
    ,@(24, $null, 16, 1, 8, 19)
    ,@(25, $null, 1, $null, $null, $null)
    ,@(26, $null, 7, 5, 9, 18)
    ,@(27, $null, 42, 22, 30, 29)
    ,@(28, $null, 60, 126, 128, 174)
    ,@(29, $null, 2, 2, 2, 5)
    ,@(30, $null, $null, $null, $null, $null)
    ,@(31, $null, 34, 35, 30, 25)
    ,@(32, $null, $null, 2, 2, 5)
    ,@(33, 2, 1, 3, 2, 3)
    ,@(34, $null, 2, 2, $null, $null)
)

$cols = @("BI", "BJ", "BK", "BL", "BM")
foreach ($entry in $rowData) {
    $r = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $val = $entry[$i + 1]
        $addr = "$($cols[$i])$r"
        if ($null -eq $val) {
            $ws.Range($addr).ClearContents()
        } else {
            $ws.Range($addr).Value2 = $val
        }
    }
}
